# Add a new credits row for the "AreaKilometer50" font used in the
# new main-menu / start-new-game UI work.
# Write the Source and Licensing columns first, then the File Name
# column, so new shared-string entries land in the same order as the
# authored workbook (Source URL, Licensing, then File Name).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "https://www.fontspace.com/a-area-kilometer-50-font-f53888"
$ws.Range("C3").Value = "Freeware"
$ws.Range("A3").Value = "AreaKilometer50-ow3xB.ttf"
# Notes/Other (D3) intentionally left blank.

# Move/refresh the active selection onto the newly added row, matching
# the saved view state of the edited workbook.
$ws.Range("D3").Select()
